$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume data per latest scrape

$ws.Range('D2').Value = '36.424.98'
$ws.Range('E2').Value = '  +0.37%  '

$ws.Range('D3').Value = '1.941.22'
$ws.Range('E3').Value = '  -1.02%  '

$ws.Range('E4').Value = '  -0.18%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '243.17'
$ws.Range('E5').Value = '  -0.09%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.612'
$ws.Range('E6').Value = '  -1.26%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '57.21'
$ws.Range('E8').Value = '  -0.10%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.360'
$ws.Range('E9').Value = '  -2.02%  '

$ws.Range('E10').Value = '  -0.72%  '

$ws.Range('E11').Value = '  -0.89%  '

$ws.Range('D12').Value = '2.227.60'

$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.811'
$ws.Range('E13').Value = '  -2.47%  '

$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '21.25'
$ws.Range('E14').Value = '  -3.17%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '13.46'
$ws.Range('E15').Value = '  -0.07%  '

$ws.Range('E16').Value = '  -3.22%  '

$ws.Range('D17').Value = '1.941.73'
$ws.Range('E17').Value = '  -1.61%  '

$ws.Range('D18').Value = '36.399.27'
$ws.Range('E18').Value = '  +0.66%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '69.28'
$ws.Range('E19').Value = '  -2.51%  '

$ws.Range('D20').Value = '0.0₃0862'
$ws.Range('E20').Value = '  -2.22%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '228.57'
$ws.Range('E21').Value = '  -0.95%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.99'
$ws.Range('E22').Value = '  -2.51%  '

$ws.Range('E23').Value = '  -0.18%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.35'
$ws.Range('E24').Value = '  -5.92%  '

$ws.Range('E25').Value = '  +0.89%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.19'
$ws.Range('E26').Value = '  -4.02%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '161.76'
$ws.Range('E27').Value = '  -2.36%  '

$ws.Range('E28').Value = '  +3.78%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.17'
$ws.Range('E29').Value = '  -3.36%  '

$ws.Range('E30').Value = '  -0.60%  '

$ws.Range('E31').Value = '  -4.82%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.57'
$ws.Range('E32').Value = '  -3.41%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0616'
$ws.Range('E33').Value = '  -3.83%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.25'
$ws.Range('E34').Value = '  +5.01%  '

$ws.Range('E35').Value = '  -4.04%  '

$ws.Range('E36').Value = '  -0.11%  '

$ws.Range('E37').Value = '  -0.80%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.16'
$ws.Range('E38').Value = '  +0.03%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.17'
$ws.Range('E39').Value = '  +9.62%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0989'
$ws.Range('E40').Value = '  +3.40%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.88'
$ws.Range('E41').Value = '  +0.23%  '

$ws.Range('E42').Value = '  -0.40%  '

$ws.Range('E43').Value = '  -3.01%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '15.93'
$ws.Range('E44').Value = '  +1.31%  '

$ws.Range('D45').Value = '1.342.00'
$ws.Range('E45').Value = '  -0.15%  '

$ws.Range('E46').Value = '  -2.85%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '86.60'
$ws.Range('E47').Value = '  -2.15%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.19'
$ws.Range('E48').Value = '  -1.04%  '

$ws.Range('E49').Value = '  +0.58%  '

$ws.Range('D50').Value = '2.118.68'

$ws.Range('E51').Value = '  -3.03%  '
